$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(4, 3).Value = 1117.911191187252
$ws.Cells.Item(5, 2).Value = '(''CASH'', ''BTC'')'
$ws.Cells.Item(5, 3).Value = 3678.390393574442
$ws.Cells.Item(6, 3).Value = 6028.238788274816
$ws.Cells.Item(7, 3).Value = 7961.788685852841
$ws.Cells.Item(8, 3).Value = 9515.697878995721
$ws.Cells.Item(9, 3).Value = 10758.29109930866
$ws.Cells.Item(10, 3).Value = 11756.3909784408
$ws.Cells.Item(11, 3).Value = 12565.80980669613
$ws.Cells.Item(12, 3).Value = 13229.93716176259
$ws.Cells.Item(13, 3).Value = 13781.51947641487
$ws.Cells.Item(14, 3).Value = 14245.05992379383
$ws.Cells.Item(15, 3).Value = 14638.93304735142
$ws.Cells.Item(16, 2).Value = '(''BTC'', ''SOL'')'
$ws.Cells.Item(16, 3).Value = 15211.89986074792
$ws.Cells.Item(17, 3).Value = 15923.26482345805
$ws.Cells.Item(18, 3).Value = 16599.10622639934
$ws.Cells.Item(19, 3).Value = 17242.70437416959
$ws.Cells.Item(20, 3).Value = 17856.75250142315
$ws.Cells.Item(21, 3).Value = 18443.50513838347
$ws.Cells.Item(22, 3).Value = 19004.88523161054
$ws.Cells.Item(23, 3).Value = 19542.56165720451
$ws.Cells.Item(24, 3).Value = 20058.00550808995
$ws.Cells.Item(25, 3).Value = 20552.53117192881
$ws.Cells.Item(26, 3).Value = 21027.32650792913
$ws.Cells.Item(27, 3).Value = 21483.47520261846
$ws.Cells.Item(28, 3).Value = 21921.97350460998
$ws.Cells.Item(29, 3).Value = 22343.74290928935
$ws.Cells.Item(30, 3).Value = 22749.63991530488
$ws.Cells.Item(31, 3).Value = 23140.46365448619
$ws.Cells.Item(32, 3).Value = 23516.96196848963
$ws.Cells.Item(33, 3).Value = 23879.83634267903
$ws.Cells.Item(34, 3).Value = 24229.74599165951
$ws.Cells.Item(35, 3).Value = 24567.3113080572
$ws.Cells.Item(36, 3).Value = 24893.11682701508
$ws.Cells.Item(37, 3).Value = 25207.71381664692
$ws.Cells.Item(38, 3).Value = 25511.62257450789
$ws.Cells.Item(39, 3).Value = 25805.33448855257
$ws.Cells.Item(40, 2).Value = '(''BTC'', ''SOL'')'
$ws.Cells.Item(40, 3).Value = 26089.3139055937
$ws.Cells.Item(41, 2).Value = '(''BTC'', ''SOL'')'
$ws.Cells.Item(41, 3).Value = 26363.99983919785
$ws.Cells.Item(42, 3).Value = 26644.533093483
$ws.Cells.Item(43, 3).Value = 43403.21323740338
